$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "Benefit" column header to "Status"
$ws.Range("D1").Value = "Status"

# Convert the Benefit amounts into a binary Status flag (1 = won/has benefit, 0 = none)
$ws.Range("D3").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D15").Value = 1

# Update the active selection on the sheet
$ws.Range("G11").Select()
